$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 1333796685911233024
$ws.Cells.Item(2,2).Value = "FirstMediaCares"
$ws.Cells.Item(2,3).Value = "@bruntus182 Selamat Malam First People. Mhn maaf sebelumnya Jika kami cek untuk area Ibu/Bapak terpantau sdh normal, mohon info saat ini kendala pd layanan Internet atau TV ya? mengenai jaringan dan tagihan Bapak / Ibu bisa mengaksesnya di situs https://t.co/h46Z2K7k7Z . Tks ^dhp"
$ws.Cells.Item(2,4).Value = "2020-12-01 15:34:45"
$ws.Cells.Item(2,5).Value = "firstmedia"
$ws.Cells.Item(2,6).Value = "0"

$ws.Cells.Item(3,1).Value = 1333794758850842880
$ws.Cells.Item(3,2).Value = "FirstMediaCares"
$ws.Cells.Item(3,3).Value = "@Ferry14939796 Hi First People, baik Pak kami akan upayakan utk perbaikannya mohon kesediaannya menunggu progres team terkait utk info selanjutnya bisa cek https://t.co/JMlireD1it terima kasih. ^rml"
$ws.Cells.Item(3,4).Value = "2020-12-01 15:27:05"
$ws.Cells.Item(3,5).Value = "firstmedia"
$ws.Cells.Item(3,6).Value = "0"

$ws.Cells.Item(4,1).Value = 1333794130921615104
$ws.Cells.Item(4,2).Value = "FirstMediaCares"
$ws.Cells.Item(4,3).Value = "@RsDayshe Selamat Malam First People, mohon maaf atas ketidaknyamanannya ya atas gangguan jaringan on-off saat ini mohon kesediaannya menunggu progres team terkait ya utk info lebih lanjut cek https://t.co/JMlireD1it terima kasih. ^rml"
$ws.Cells.Item(4,4).Value = "2020-12-01 15:24:36"
$ws.Cells.Item(4,5).Value = "firstmedia"
$ws.Cells.Item(4,6).Value = "0"

$ws.Cells.Item(5,1).Value = 1333790058327855104
$ws.Cells.Item(5,2).Value = "FirstMediaCares"
$ws.Cells.Item(5,3).Value = "@godjila526 Selamat malam first people,mhn maaf atas ketidaknyamanan.Saat ini area anda mengalami gangguan quality signal sehingga internet slow, on off bahkan offline dan saat ini masih dalam proses perbaikan oleh tim terkait.Detailnya https://t.co/h46Z2K7k7Z.Tks ^Fjr"
$ws.Cells.Item(5,4).Value = "2020-12-01 15:08:25"
$ws.Cells.Item(5,5).Value = "firstmedia"
$ws.Cells.Item(5,6).Value = "0"

$ws.Cells.Item(6,1).Value = 1333787153592058112
$ws.Cells.Item(6,2).Value = "FirstMediaCares"
$ws.Cells.Item(6,3).Value = "@nina_mayoni Selamat malam first people,mhn maaf atas ketidaknyamanan.Saat ini area anda mengalami gangguan quality signal sehingga internet on off bahkan offline dan saat ini masih dalam proses perbaikan oleh tim terkait.Detailnya https://t.co/h46Z2K7k7Z.Tks ^Fjr"
$ws.Cells.Item(6,4).Value = "2020-12-01 14:56:52"
$ws.Cells.Item(6,5).Value = "firstmedia"
$ws.Cells.Item(6,6).Value = "0"

$ws.Cells.Item(7,1).Value = 1333782113636684032
$ws.Cells.Item(7,2).Value = "alphavityazi"
$ws.Cells.Item(7,3).Value = "@FirstMediaCares setelah dicek, lalu apa? ini udah sekitar 2-3 hari gangguan di website cek firstmedia, tp ga ada perubahan. gausah perubahan, estimasi perbaikan aja gak jelas s/d kapan. belum lagi rto rto yang terjadi sudah hampir 2 minggu lebih."
$ws.Cells.Item(7,4).Value = "2020-12-01 14:36:51"
$ws.Cells.Item(7,5).Value = "firstmedia"
$ws.Cells.Item(7,6).Value = "0"

$ws.Cells.Item(8,1).Value = 1333769726556705024
$ws.Cells.Item(8,2).Value = "FirstMediaCares"
$ws.Cells.Item(8,3).Value = "@redblueexh Hi first people, mohon maaf atas ketidaknyamanannya perihal gangguan jaringan di area rmh Bapak/Ibu saat ini dan msh dalam proses perbaikan. Update infonya ada di https://t.co/h46Z2K7k7Z ya. Tks ^tst"
$ws.Cells.Item(8,4).Value = "2020-12-01 13:47:37"
$ws.Cells.Item(8,5).Value = "firstmedia"
$ws.Cells.Item(8,6).Value = "0"

$ws.Cells.Item(9,1).Value = 1333767129007493120
$ws.Cells.Item(9,2).Value = "FirstMediaCares"
$ws.Cells.Item(9,3).Value = "@OgawaRen Selamat malam First People, mohon maaf atas ketidaknyamanannya. Saat ini area Bpk/Ibu sdg mengalami gangguan jaringan. Namun akan kami upayakan secepatnya. Informasi mengenai status area tsb silahkan dicek di https://t.co/IQIwB9XHrI. Terimakasih ^ang"
$ws.Cells.Item(9,4).Value = "2020-12-01 13:37:18"
$ws.Cells.Item(9,5).Value = "firstmedia"
$ws.Cells.Item(9,6).Value = "0"

$ws.Cells.Item(10,1).Value = 1333760783608528896
$ws.Cells.Item(10,2).Value = "fuckinghard6969"
$ws.Cells.Item(10,3).Value = "@febriani_p @FirstMediaCares Emang aneh2 aja first media skrg.. boikot aja firstmedia"
$ws.Cells.Item(10,4).Value = "2020-12-01 13:12:05"
$ws.Cells.Item(10,5).Value = "firstmedia"
$ws.Cells.Item(10,6).Value = "0"

$ws.Cells.Item(11,1).Value = 1333758666256708096
$ws.Cells.Item(11,2).Value = "FirstMediaCares"
$ws.Cells.Item(11,3).Value = "@davfadel Selamat Malam First People. Mohon maaf saat ini sdg terjadi prblm jaringan dan hal ini sdh dalam proses perbaikan dari team jaringan kami. Estimasi waktu perbaikan 02-Dec-2020 / 07:58. Detailnya bisa pantau di https://t.co/h46Z2K7k7Z. Tks ^Ibn"
$ws.Cells.Item(11,4).Value = "2020-12-01 13:03:40"
$ws.Cells.Item(11,5).Value = "firstmedia"
$ws.Cells.Item(11,6).Value = "0"

$ws.Cells.Item(12,1).Value = 1333740413656527104
$ws.Cells.Item(12,2).Value = "FirstMediaCares"
$ws.Cells.Item(12,3).Value = "@Aldry_F Selamat malam First People. Mhn maaf atas ketidaknyamanannya. Kami cek area saat ini normal, mhn restart kabel pwr modem ya. Jika masih kendala mhn dicoba utk bypass dr modem lgsg dgn kabel LAN ya. Untuk cek &amp; refresh jaringan bs di https://t.co/h46Z2K7k7Z ya. Tks. ^Fjr"
$ws.Cells.Item(12,4).Value = "2020-12-01 11:51:08"
$ws.Cells.Item(12,5).Value = "firstmedia"
$ws.Cells.Item(12,6).Value = "0"

$ws.Cells.Item(13,1).Value = 1333739474392470016
$ws.Cells.Item(13,2).Value = "FirstMediaCares"
$ws.Cells.Item(13,3).Value = "@sen9922 Selamat malam FIRST people. Mohon maaf atas ketidaknyamanan yang dialami. Saat ini areanya masih dalam prses perbaikan, yang menyebabkan kualitas jaringan areannya  menurun.  untuk cek status perbaikan bisa cek di https://t.co/h46Z2K7k7Z . Tks ^Fuj"
$ws.Cells.Item(13,4).Value = "2020-12-01 11:47:25"
$ws.Cells.Item(13,5).Value = "firstmedia"
$ws.Cells.Item(13,6).Value = "0"

$ws.Cells.Item(14,1).Value = 1333672686745121024
$ws.Cells.Item(14,2).Value = "FirstMediaCares"
$ws.Cells.Item(14,3).Value = "@rizkybayumilano   Selamat Siang First People. Mohon maaf atas ketidaknyamanannya. Saat ini sdg terjadi prblm jaringan signal quality. Detail jaringan dan tagihan bisa pantau di https://t.co/h46Z2K7k7Z  dan MYFirstMedia App ^em"
$ws.Cells.Item(14,4).Value = "2020-12-01 07:22:01"
$ws.Cells.Item(14,5).Value = "firstmedia"
$ws.Cells.Item(14,6).Value = "0"

$ws.Cells.Item(15,1).Value = 1333654353995973120
$ws.Cells.Item(15,2).Value = "FirstMediaCares"
$ws.Cells.Item(15,3).Value = "@rahmat_sitinjak Slmt Siang First People. Mhn maaf atas ktdknyamanannya. Sdg terjadi gangguan layanan area ASEM BARIS, adapun perkiraan estimasi perbaikan sampai dgn 01-Dec-2020 / 21:18. Utk info terakhir mengenai mslh ini dpt mengaksesnya via https://t.co/h46Z2K7k7Z. Tks ^mrs"
$ws.Cells.Item(15,4).Value = "2020-12-01 06:09:10"
$ws.Cells.Item(15,5).Value = "firstmedia"
$ws.Cells.Item(15,6).Value = "0"

$ws.Cells.Item(16,1).Value = 1333649804295959040
$ws.Cells.Item(16,2).Value = "FirstMediaCares"
$ws.Cells.Item(16,3).Value = "@soreschach Selamat Siang First Poeple. Mohon maaf atas ketidaknyamanannya. Saat ini sdg ada gangguan jaringan di area Bapak/Ibu dan msh dlm proses perbaikan oleh tim terkait kami. Untuk updatenya bisa akses https://t.co/h46Z2K7k7Z. Tks ^RZA"
$ws.Cells.Item(16,4).Value = "2020-12-01 05:51:06"
$ws.Cells.Item(16,5).Value = "firstmedia"
$ws.Cells.Item(16,6).Value = "0"

$ws.Cells.Item(17,1).Value = 1333644522337693952
$ws.Cells.Item(17,2).Value = "FirstMediaCares"
$ws.Cells.Item(17,3).Value = "@boreg8401 Hi First People. Mhn maaf atas kendalanya dan terima kasih inputannya untuk perbaikan layanan kami. Untuk update info status jaringan kami sarankan untuk cek di https://t.co/h46Z2K7k7Z ya. Tks ^mel"
$ws.Cells.Item(17,4).Value = "2020-12-01 05:30:06"
$ws.Cells.Item(17,5).Value = "firstmedia"
$ws.Cells.Item(17,6).Value = "0"

$ws.Cells.Item(18,1).Value = 1333624515889221888
$ws.Cells.Item(18,2).Value = "FirstMediaCares"
$ws.Cells.Item(18,3).Value = "@kabigon89    Selamat  pagi  First People. Mohon maaf atas ketidaknyamanannya. Saat ini sdg terjadi gangguan jaringan. estimasi waktu 01-Dec-2020 / 22:52 kami sedang upayakan segera normal,  Detailnya bisa pantau di https://t.co/h46Z2K7k7Z. Tks ^Sth"
$ws.Cells.Item(18,4).Value = "2020-12-01 04:10:36"
$ws.Cells.Item(18,5).Value = "firstmedia"
$ws.Cells.Item(18,6).Value = "0"

$ws.Cells.Item(19,1).Value = 1333614841504796928
$ws.Cells.Item(19,2).Value = "rezapati"
$ws.Cells.Item(19,3).Value = "@FirstMediaCares parah nih firstmedia dalam seminggu 2x gangguan. Perusahaan besar dan lama tapi belum bisa meminimalisir gangguan. Ditambah disaat orang sedang butuh koneksi internet. Mmm.. sudah waktunya cari-cari alternatif lain nih."
$ws.Cells.Item(19,4).Value = "2020-12-01 03:32:10"
$ws.Cells.Item(19,5).Value = "firstmedia"
$ws.Cells.Item(19,6).Value = "0"

$ws.Cells.Item(20,1).Value = 1333590577875623936
$ws.Cells.Item(20,2).Value = "notapetite"
$ws.Cells.Item(20,3).Value = "firstmedia masi pagi udh gajelas☺️🙏🏻"
$ws.Cells.Item(20,4).Value = "2020-12-01 01:55:45"
$ws.Cells.Item(20,5).Value = "firstmedia"
$ws.Cells.Item(20,6).Value = "0"

$ws.Cells.Item(21,2).Value = "FirstMediaCares"
$ws.Cells.Item(21,3).Value = "@HafidRehaff Selamat Pagi First People, mohon maaf atas ketidaknyamanannya. Kami informasikan untuk saat ini area sedang mengalami gangguan, estimasi waktu perbaikan diperkirakan 02-Dec-2020 / 01:26. Info status jaringan dpt akses di https://t.co/h46Z2K7k7Z. Terima kasih ^nlv"
$ws.Cells.Item(21,4).Value = "2020-12-01 01:54:51"
$ws.Cells.Item(21,5).Value = "firstmedia"
$ws.Cells.Item(21,6).Value = "0"

